$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 277
$ws.Range("I8").Value = 277
$ws.Range("K8").Value = 831
$ws.Range("M8").Value = -692

# Row 9
$ws.Range("H9").Value = 176.06667
$ws.Range("I9").Value = 127.888885
$ws.Range("J9").Value = 248.33333
$ws.Range("K9").Value = 127.888885
$ws.Range("L9").Value = 248.33333
$ws.Range("M9").Value = 41.111115
$ws.Range("N9").Value = -586.3333299999999

# Row 41
$ws.Range("H41").Value = 646.1875
$ws.Range("I41").Value = 361.58334
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 361.58334
$ws.Range("L41").Value = 1500
$ws.Range("M41").Value = 78.41665999999998
$ws.Range("N41").Value = -2380

# Row 76
$ws.Range("H76").Value = 3730.6943
$ws.Range("J76").Value = 4230
$ws.Range("L76").Value = 4230
$ws.Range("N76").Value = -4860

# Row 79
$ws.Range("H79").Value = 3730.6943
$ws.Range("J79").Value = 4230
$ws.Range("L79").Value = 4230
$ws.Range("N79").Value = -6414

# Row 87
$ws.Range("H87").Value = 35293.582
$ws.Range("J87").Value = 35293.582
$ws.Range("L87").Value = 35293.582
$ws.Range("N87").Value = -37789.582

# Row 90
$ws.Range("H90").Value = 35293.582
$ws.Range("J90").Value = 35293.582
$ws.Range("L90").Value = 105880.746
$ws.Range("N90").Value = -118360.746

# Row 112
$ws.Range("H112").Value = 1239.8644
$ws.Range("J112").Value = 1239.8644
$ws.Range("L112").Value = 3719.5932
$ws.Range("N112").Value = -5935.593199999999

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 138
$ws.Range("H138").Value = 3702.1646
$ws.Range("I138").Value = 2497.375
$ws.Range("J138").Value = 4008.1428
$ws.Range("K138").Value = 7492.125
$ws.Range("L138").Value = 12024.4284
$ws.Range("M138").Value = -2352.125
$ws.Range("N138").Value = -22304.4284

$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 61
$ws.Range("H61").Value = 5034.2354
$ws.Range("I61").Value = 3863.6956
$ws.Range("J61").Value = 7481.727
$ws.Range("K61").Value = 3863.6956
$ws.Range("L61").Value = 7481.727
$ws.Range("M61").Value = -3651.6956
$ws.Range("N61").Value = -7905.727

# Row 74
$ws.Range("H74").Value = 5975.355
$ws.Range("I74").Value = 3301.1428
$ws.Range("J74").Value = 11591.2
$ws.Range("K74").Value = 3301.1428
$ws.Range("L74").Value = 11591.2
$ws.Range("M74").Value = -2427.1428
$ws.Range("N74").Value = -13339.2

# Row 77
$ws.Range("H77").Value = 5975.355
$ws.Range("I77").Value = 3301.1428
$ws.Range("J77").Value = 11591.2
$ws.Range("K77").Value = 16505.714
$ws.Range("L77").Value = 57956
$ws.Range("M77").Value = -12137.714
$ws.Range("N77").Value = -66692

# Row 136
$ws.Range("H136").Value = 5034.2354
$ws.Range("I136").Value = 3863.6956
$ws.Range("J136").Value = 7481.727
$ws.Range("K136").Value = 11591.0868
$ws.Range("L136").Value = 22445.181
$ws.Range("M136").Value = -9041.086800000001
$ws.Range("N136").Value = -27545.181

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3960.5
$ws.Range("I134").Value = 3847.1943
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 11541.5829
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -9006.582900000001
$ws.Range("N134").Value = -23070

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 83.86667
$ws.Range("I7").Value = 81
$ws.Range("J7").Value = 88.166664
$ws.Range("K7").Value = 81
$ws.Range("L7").Value = 88.166664
$ws.Range("M7").Value = 32
$ws.Range("N7").Value = -314.166664

# Row 31
$ws.Range("H31").Value = 1496.0618
$ws.Range("I31").Value = 1135.6528
$ws.Range("J31").Value = 4379.3335
$ws.Range("K31").Value = 1135.6528
$ws.Range("L31").Value = 4379.3335
$ws.Range("M31").Value = -840.6528000000001
$ws.Range("N31").Value = -4969.3335

# Row 34
$ws.Range("H34").Value = 1496.0618
$ws.Range("I34").Value = 1135.6528
$ws.Range("J34").Value = 4379.3335
$ws.Range("K34").Value = 1135.6528
$ws.Range("L34").Value = 4379.3335
$ws.Range("M34").Value = -933.6528000000001
$ws.Range("N34").Value = -4783.3335

# Row 86
$ws.Range("H86").Value = 1930.8
$ws.Range("I86").Value = 1600
$ws.Range("J86").Value = 2151.3333
$ws.Range("K86").Value = 1600
$ws.Range("L86").Value = 2151.3333
$ws.Range("M86").Value = -477
$ws.Range("N86").Value = -4397.3333

# Row 89
$ws.Range("H89").Value = 1930.8
$ws.Range("I89").Value = 1600
$ws.Range("J89").Value = 2151.3333
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 10756.6665
$ws.Range("M89").Value = -2384
$ws.Range("N89").Value = -21988.6665

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 122
$ws.Range("H122").Value = 11888.8125
$ws.Range("I122").Value = 4015.923
$ws.Range("J122").Value = 46004.668
$ws.Range("K122").Value = 12047.769
$ws.Range("L122").Value = 138014.004
$ws.Range("M122").Value = -9597.769
$ws.Range("N122").Value = -142914.004

# Row 132
$ws.Range("H132").Value = 2632.087
$ws.Range("I132").Value = 2238.9473
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 6716.841899999999
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -4186.841899999999
$ws.Range("N132").Value = -18558.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 16670547
$ws.Range("I5").Value = 414.66666
$ws.Range("J5").Value = 83351080
$ws.Range("K5").Value = 1243.99998
$ws.Range("L5").Value = 250053240
$ws.Range("M5").Value = -1131.99998
$ws.Range("N5").Value = -250053464

# Row 107
$ws.Range("H107").Value = 1155
$ws.Range("I107").Value = 334.22223
$ws.Range("J107").Value = 1360.1945
$ws.Range("K107").Value = 1002.66669
$ws.Range("L107").Value = 4080.5835
$ws.Range("M107").Value = 917.33331
$ws.Range("N107").Value = -7920.583500000001

# Row 120
$ws.Range("H120").Value = 10220.19
$ws.Range("I120").Value = 4984
$ws.Range("J120").Value = 11856.5
$ws.Range("K120").Value = 14952
$ws.Range("L120").Value = 35569.5
$ws.Range("M120").Value = -10114
$ws.Range("N120").Value = -45245.5

# Row 131
$ws.Range("H131").Value = 815.2
$ws.Range("I131").Value = 475
$ws.Range("J131").Value = 848.8461
$ws.Range("K131").Value = 1425
$ws.Range("L131").Value = 2546.5383
$ws.Range("M131").Value = 3615
$ws.Range("N131").Value = -12626.5383

# Row 135
$ws.Range("H135").Value = 16670547
$ws.Range("I135").Value = 414.66666
$ws.Range("J135").Value = 83351080
$ws.Range("K135").Value = 3731.99994
$ws.Range("L135").Value = 750159720
$ws.Range("M135").Value = -1196.99994
$ws.Range("N135").Value = -750164790

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

# Row 80
$ws.Range("H80").Value = 6063.125
$ws.Range("I80").Value = 8637.25
$ws.Range("K80").Value = 8637.25
$ws.Range("M80").Value = -7639.25

# Row 83
$ws.Range("H83").Value = 6063.125
$ws.Range("I83").Value = 8637.25
$ws.Range("K83").Value = 43186.25
$ws.Range("M83").Value = -38194.25

# Row 132
$ws.Range("H132").Value = 3819.3333
$ws.Range("I132").Value = 2659.6667
$ws.Range("K132").Value = 7979.000100000001
$ws.Range("M132").Value = -5449.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1683.3334
$ws.Range("I68").Value = 1575
$ws.Range("K68").Value = 1575
$ws.Range("M68").Value = -826

# Row 71
$ws.Range("H71").Value = 1683.3334
$ws.Range("I71").Value = 1575
$ws.Range("K71").Value = 7875
$ws.Range("M71").Value = -4131

# Row 136
$ws.Range("H136").Value = 6063.0356
$ws.Range("I136").Value = 5015.475
$ws.Range("K136").Value = 15046.425
$ws.Range("M136").Value = -12496.425

$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 26
$ws.Range("H26").Value = 8100
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 8100
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 8100
$ws.Range("N26").Value = -8686
$ws.Range("M26").ClearContents()

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 136
$ws.Range("H136").Value = 5080.94
$ws.Range("I136").Value = 2403.724
$ws.Range("J136").Value = 8778.048000000001
$ws.Range("K136").Value = 7211.172
$ws.Range("L136").Value = 26334.144
$ws.Range("M136").Value = -4661.172
$ws.Range("N136").Value = -31434.144
